$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 766
$ws.Range("F3").Value = 14304
$ws.Range("F4").Value = 14432
$ws.Range("F5").Value = 1366
$ws.Range("F6").Value = 1406
$ws.Range("F7").Value = 5920
$ws.Range("F8").Value = 987
$ws.Range("F9").Value = 578
$ws.Range("F12").Value = 200
$ws.Range("F13").Value = 1561
$ws.Range("F14").Value = 448
$ws.Range("F15").Value = 2132
$ws.Range("F16").Value = 1213
$ws.Range("F17").Value = 1855
$ws.Range("F18").Value = 918
$ws.Range("F19").Value = 36
$ws.Range("F21").Value = 571
$ws.Range("F22").Value = 822
$ws.Range("F23").Value = 3356
$ws.Range("F25").Value = 315
$ws.Range("F26").Value = 2425
$ws.Range("F27").Value = 603
$ws.Range("F29").Value = 1340
$ws.Range("F30").Value = 1805
$ws.Range("F31").Value = 1077
$ws.Range("F32").Value = 1418
$ws.Range("F33").Value = 108
$ws.Range("F35").Value = 4902
$ws.Range("F36").Value = 4893
$ws.Range("F37").Value = 309
$ws.Range("F38").Value = 162
$ws.Range("F39").Value = 680
$ws.Range("F40").Value = 692
$ws.Range("F41").Value = 3308
$ws.Range("F43").Value = 925
$ws.Range("F44").Value = 343
$ws.Range("F45").Value = 113
$ws.Range("F46").Value = 94
$ws.Range("F47").Value = 4439
$ws.Range("F48").Value = 608
$ws.Range("F49").Value = 298
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 122
$ws.Range("F20").Value = 15
$ws.Range("F21").Value = 109
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7627
$ws.Range("F3").Value = 249
$ws.Range("F4").Value = 839
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 766
$ws.Range("F3").Value = 249
$ws.Range("F4").Value = 839
$ws.Range("F6").Value = 14305
$ws.Range("F8").Value = 1406
$ws.Range("F9").Value = 5920
$ws.Range("F10").Value = 987
$ws.Range("F11").Value = 122
$ws.Range("F14").Value = 1561
$ws.Range("F15").Value = 448
$ws.Range("F16").Value = 1213
$ws.Range("F17").Value = 1855
$ws.Range("F18").Value = 918
$ws.Range("F19").Value = 36
$ws.Range("F20").Value = 571
$ws.Range("F21").Value = 3356
$ws.Range("F22").Value = 315
$ws.Range("F23").Value = 603
$ws.Range("F25").Value = 1805
$ws.Range("F28").Value = 1418
$ws.Range("F30").Value = 108
$ws.Range("F32").Value = 4902
$ws.Range("F33").Value = 4893
$ws.Range("F35").Value = 309
$ws.Range("F36").Value = 162
$ws.Range("F37").Value = 680
$ws.Range("F38").Value = 692
$ws.Range("F39").Value = 3308
$ws.Range("F41").Value = 925
$ws.Range("F42").Value = 343
$ws.Range("F44").Value = 94
$ws.Range("F45").Value = 4439
$ws.Range("F46").Value = 608
$ws.Range("F47").Value = 298
